$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 3381.6155  # ALC H19 was 3389.8462
$ws.Cells.Item(19, 9).Value = 3027.5  # ALC I19 was 2761.4
$ws.Cells.Item(19, 10).Value = 3539  # ALC J19 was 3782.625
$ws.Cells.Item(19, 11).Value = 3027.5  # ALC K19 was 2761.4
$ws.Cells.Item(19, 12).Value = 3539  # ALC L19 was 3782.625
$ws.Cells.Item(19, 13).Value = -2852.5  # ALC M19 was -2586.4
$ws.Cells.Item(19, 14).Value = -3889  # ALC N19 was -4132.625
$ws.Cells.Item(33, 8).Value = 16991.842  # ALC H33 was 15399.762
$ws.Cells.Item(33, 9).Value = 17880.277  # ALC I33 was 16941.842
$ws.Cells.Item(33, 10).Value = 1000  # ALC J33 was 750
$ws.Cells.Item(33, 11).Value = 17880.277  # ALC K33 was 16941.842
$ws.Cells.Item(33, 12).Value = 1000  # ALC L33 was 750
$ws.Cells.Item(33, 13).Value = -17651.277  # ALC M33 was -16712.842
$ws.Cells.Item(33, 14).Value = -1458  # ALC N33 was -1208
$ws.Cells.Item(80, 8).Value = 311.7143  # ALC H80 was 276.2143
$ws.Cells.Item(80, 9).Value = 263.75  # ALC I80 was 251.3077
$ws.Cells.Item(80, 10).Value = 599.5  # ALC J80 was 600
$ws.Cells.Item(80, 11).Value = 791.25  # ALC K80 was 753.9231
$ws.Cells.Item(80, 12).Value = 1798.5  # ALC L80 was 1800
$ws.Cells.Item(80, 13).Value = 206.75  # ALC M80 was 244.0769
$ws.Cells.Item(80, 14).Value = -3794.5  # ALC N80 was -3796
$ws.Cells.Item(83, 8).Value = 311.7143  # ALC H83 was 276.2143
$ws.Cells.Item(83, 9).Value = 263.75  # ALC I83 was 251.3077
$ws.Cells.Item(83, 10).Value = 599.5  # ALC J83 was 600
$ws.Cells.Item(83, 11).Value = 2373.75  # ALC K83 was 2261.7693
$ws.Cells.Item(83, 12).Value = 5395.5  # ALC L83 was 5400
$ws.Cells.Item(83, 13).Value = 2618.25  # ALC M83 was 2730.2307
$ws.Cells.Item(83, 14).Value = -15379.5  # ALC N83 was -15384
$ws.Cells.Item(94, 8).Value = 504.75  # ALC H94 was 540
$ws.Cells.Item(94, 9).Value = 504.75  # ALC I94 was 540
$ws.Cells.Item(94, 11).Value = 504.75  # ALC K94 was 540
$ws.Cells.Item(94, 13).Value = -53.75  # ALC M94 was -89
$ws.Cells.Item(132, 8).Value = 1792.4762  # ALC H132 was 1834.65
$ws.Cells.Item(132, 9).Value = 1798.2632  # ALC I132 was 1845.4445
$ws.Cells.Item(132, 11).Value = 5394.7896  # ALC K132 was 5536.333500000001
$ws.Cells.Item(132, 13).Value = -2864.7896  # ALC M132 was -3006.333500000001
$ws.Cells.Item(137, 8).Value = 1982.9333  # ALC H137 was 1495.5385
$ws.Cells.Item(137, 9).Value = 1221.5555  # ALC I137 was 1336.5
$ws.Cells.Item(137, 10).Value = 3125  # ALC J137 was 1750
$ws.Cells.Item(137, 11).Value = 3664.6665  # ALC K137 was 4009.5
$ws.Cells.Item(137, 12).Value = 9375  # ALC L137 was 5250
$ws.Cells.Item(137, 13).Value = -1114.6665  # ALC M137 was -1459.5
$ws.Cells.Item(137, 14).Value = -14475  # ALC N137 was -10350
$ws.Cells.Item(138, 8).Value = 18184016  # ALC H138 was 2206.375
$ws.Cells.Item(138, 9).Value = 1914.3  # ALC I138 was 2031.25
$ws.Cells.Item(138, 10).Value = 22224482  # ALC J138 was 2235.5625
$ws.Cells.Item(138, 11).Value = 5742.9  # ALC K138 was 6093.75
$ws.Cells.Item(138, 12).Value = 66673446  # ALC L138 was 6706.6875
$ws.Cells.Item(138, 13).Value = -602.8999999999996  # ALC M138 was -953.75
$ws.Cells.Item(138, 14).Value = -66683726  # ALC N138 was -16986.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7576.4053  # ARM H32 was 7576.4595
$ws.Cells.Item(32, 9).Value = 7437.914  # ARM I32 was 7648.794
$ws.Cells.Item(32, 10).Value = 10000  # ARM J32 was 6756.6665
$ws.Cells.Item(32, 11).Value = 7437.914  # ARM K32 was 7648.794
$ws.Cells.Item(32, 12).Value = 10000  # ARM L32 was 6756.6665
$ws.Cells.Item(32, 13).Value = -7150.914  # ARM M32 was -7361.794
$ws.Cells.Item(32, 14).Value = -10574  # ARM N32 was -7330.6665
$ws.Cells.Item(44, 8).Value = 31300  # ARM H44 was 31000
$ws.Cells.Item(44, 10).Value = 31300  # ARM J44 was 31000
$ws.Cells.Item(44, 12).Value = 31300  # ARM L44 was 31000
$ws.Cells.Item(44, 14).Value = -32276  # ARM N44 was -31976
$ws.Cells.Item(45, 8).Value = 229055.56  # ARM H45 was 256875
$ws.Cells.Item(45, 10).Value = 8700  # ARM J45 was 9250
$ws.Cells.Item(45, 12).Value = 8700  # ARM L45 was 9250
$ws.Cells.Item(45, 14).Value = -9454  # ARM N45 was -10004
$ws.Cells.Item(49, 8).Value = 0  # ARM H49 was 20000
$ws.Cells.Item(49, 9).Value = 0  # ARM I49 was 20000
$ws.Cells.Item(49, 11).Value = 0  # ARM K49 was 20000
$ws.Cells.Item(49, 13).ClearContents()  # ARM M49 was -19740

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 126499.875  # CRP H4 was 112499.89
$ws.Cells.Item(4, 10).Value = 1599.8  # CRP J4 was 1416.5
$ws.Cells.Item(4, 12).Value = 1599.8  # CRP L4 was 1416.5
$ws.Cells.Item(4, 14).Value = -1823.8  # CRP N4 was -1640.5
$ws.Cells.Item(8, 8).Value = 4966.5  # CRP H8 was 6749.75
$ws.Cells.Item(8, 9).Value = 1400  # CRP I8 was 0
$ws.Cells.Item(8, 11).Value = 1400  # CRP K8 was 0
$ws.Cells.Item(8, 13).Value = -1260  # CRP M8 new cell
$ws.Cells.Item(15, 8).Value = 5310.8887  # CRP H15 was 4710.4443
$ws.Cells.Item(15, 9).Value = 1720  # CRP I15 was 1599.1666
$ws.Cells.Item(15, 10).Value = 9799.5  # CRP J15 was 10933
$ws.Cells.Item(15, 11).Value = 1720  # CRP K15 was 1599.1666
$ws.Cells.Item(15, 12).Value = 9799.5  # CRP L15 was 10933
$ws.Cells.Item(15, 13).Value = -1550  # CRP M15 was -1429.1666
$ws.Cells.Item(15, 14).Value = -10139.5  # CRP N15 was -11273
$ws.Cells.Item(16, 8).Value = 2856.963  # CRP H16 was 2586.4666
$ws.Cells.Item(16, 9).Value = 2628.3333  # CRP I16 was 2564.2727
$ws.Cells.Item(16, 10).Value = 3039.8667  # CRP J16 was 2599.3157
$ws.Cells.Item(16, 11).Value = 2628.3333  # CRP K16 was 2564.2727
$ws.Cells.Item(16, 12).Value = 3039.8667  # CRP L16 was 2599.3157
$ws.Cells.Item(16, 13).Value = -2341.3333  # CRP M16 was -2277.2727
$ws.Cells.Item(16, 14).Value = -3613.8667  # CRP N16 was -3173.3157
$ws.Cells.Item(31, 8).Value = 8249.5  # CRP H31 was 9099.700000000001
$ws.Cells.Item(31, 9).Value = 6555  # CRP I31 was 7285.4287
$ws.Cells.Item(31, 11).Value = 6555  # CRP K31 was 7285.4287
$ws.Cells.Item(31, 13).Value = -6260  # CRP M31 was -6990.4287
$ws.Cells.Item(34, 8).Value = 8249.5  # CRP H34 was 9099.700000000001
$ws.Cells.Item(34, 9).Value = 6555  # CRP I34 was 7285.4287
$ws.Cells.Item(34, 11).Value = 6555  # CRP K34 was 7285.4287
$ws.Cells.Item(34, 13).Value = -6353  # CRP M34 was -7083.4287
$ws.Cells.Item(99, 8).Value = 5750  # CRP H99 was 5533.3335
$ws.Cells.Item(99, 10).Value = 0  # CRP J99 was 5100
$ws.Cells.Item(99, 12).Value = 0  # CRP L99 was 5100
$ws.Cells.Item(99, 14).ClearContents()  # CRP N99 was -8096
$ws.Cells.Item(113, 8).Value = 2856.963  # CRP H113 was 2586.4666
$ws.Cells.Item(113, 9).Value = 2628.3333  # CRP I113 was 2564.2727
$ws.Cells.Item(113, 10).Value = 3039.8667  # CRP J113 was 2599.3157
$ws.Cells.Item(113, 11).Value = 2628.3333  # CRP K113 was 2564.2727
$ws.Cells.Item(113, 12).Value = 3039.8667  # CRP L113 was 2599.3157
$ws.Cells.Item(113, 13).Value = -458.3332999999998  # CRP M113 was -394.2727
$ws.Cells.Item(113, 14).Value = -7379.8667  # CRP N113 was -6939.3157
$ws.Cells.Item(126, 8).Value = 5750  # CRP H126 was 5533.3335
$ws.Cells.Item(126, 10).Value = 0  # CRP J126 was 5100
$ws.Cells.Item(126, 12).Value = 0  # CRP L126 was 15300
$ws.Cells.Item(126, 14).ClearContents()  # CRP N126 was -20240

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 3100  # CUL H12 was 532.8889
$ws.Cells.Item(12, 9).Value = 0  # CUL I12 was 227.75
$ws.Cells.Item(12, 10).Value = 3100  # CUL J12 was 620.0714
$ws.Cells.Item(12, 11).Value = 0  # CUL K12 was 683.25
$ws.Cells.Item(12, 12).Value = 9300  # CUL L12 was 1860.2142
$ws.Cells.Item(12, 13).ClearContents()  # CUL M12 was -510.25
$ws.Cells.Item(12, 14).Value = -9646  # CUL N12 was -2206.2142
$ws.Cells.Item(131, 8).Value = 20001958  # CUL H131 was 22729396
$ws.Cells.Item(131, 9).Value = 83334056  # CUL I131 was 125000670
$ws.Cells.Item(131, 10).Value = 2347.2104  # CUL J131 was 2446.2222
$ws.Cells.Item(131, 11).Value = 250002168  # CUL K131 was 375002010
$ws.Cells.Item(131, 12).Value = 7041.6312  # CUL L131 was 7338.6666
$ws.Cells.Item(131, 13).Value = -249997128  # CUL M131 was -374996970
$ws.Cells.Item(131, 14).Value = -17121.6312  # CUL N131 was -17418.6666
$ws.Cells.Item(132, 8).Value = 41668148  # CUL H132 was 1688.8
$ws.Cells.Item(132, 9).Value = 62501224  # CUL I132 was 1481.3334
$ws.Cells.Item(132, 11).Value = 562511016  # CUL K132 was 13332.0006
$ws.Cells.Item(132, 13).Value = -562508486  # CUL M132 was -10802.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 20000  # GSM H21 was 14000
$ws.Cells.Item(21, 10).Value = 20000  # GSM J21 was 14000
$ws.Cells.Item(21, 12).Value = 20000  # GSM L21 was 14000
$ws.Cells.Item(21, 14).Value = -20346  # GSM N21 was -14346
$ws.Cells.Item(30, 8).Value = 20000  # GSM H30 was 14000
$ws.Cells.Item(30, 10).Value = 20000  # GSM J30 was 14000
$ws.Cells.Item(30, 12).Value = 20000  # GSM L30 was 14000
$ws.Cells.Item(30, 14).Value = -20210  # GSM N30 was -14210
$ws.Cells.Item(122, 8).Value = 3629.2856  # GSM H122 was 3629.4285
$ws.Cells.Item(122, 9).Value = 2802  # GSM I122 was 2802.3333
$ws.Cells.Item(122, 11).Value = 8406  # GSM K122 was 8406.999899999999
$ws.Cells.Item(122, 13).Value = -5956  # GSM M122 was -5956.999899999999
$ws.Cells.Item(126, 8).Value = 4298.8  # GSM H126 was 3617.4285
$ws.Cells.Item(126, 9).Value = 3166  # GSM I126 was 2564.4
$ws.Cells.Item(126, 10).Value = 5998  # GSM J126 was 6250
$ws.Cells.Item(126, 11).Value = 9498  # GSM K126 was 7693.200000000001
$ws.Cells.Item(126, 12).Value = 17994  # GSM L126 was 18750
$ws.Cells.Item(126, 13).Value = -7028  # GSM M126 was -5223.200000000001
$ws.Cells.Item(126, 14).Value = -22934  # GSM N126 was -23690

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3724.75  # LTW H22 was 3179.8
$ws.Cells.Item(22, 9).Value = 3724.75  # LTW I22 was 3179.8
$ws.Cells.Item(22, 11).Value = 3724.75  # LTW K22 was 3179.8
$ws.Cells.Item(22, 13).Value = -3429.75  # LTW M22 was -2884.8
$ws.Cells.Item(27, 8).Value = 3724.75  # LTW H27 was 3179.8
$ws.Cells.Item(27, 9).Value = 3724.75  # LTW I27 was 3179.8
$ws.Cells.Item(27, 11).Value = 3724.75  # LTW K27 was 3179.8
$ws.Cells.Item(27, 13).Value = -3617.75  # LTW M27 was -3072.8
$ws.Cells.Item(46, 8).Value = 16412.824  # LTW H46 was 15578.579
$ws.Cells.Item(46, 9).Value = 3024.75  # LTW I46 was 2814.8
$ws.Cells.Item(46, 10).Value = 20532.23  # LTW J46 was 20137.072
$ws.Cells.Item(46, 11).Value = 3024.75  # LTW K46 was 2814.8
$ws.Cells.Item(46, 12).Value = 20532.23  # LTW L46 was 20137.072
$ws.Cells.Item(46, 13).Value = -2836.75  # LTW M46 was -2626.8
$ws.Cells.Item(46, 14).Value = -20908.23  # LTW N46 was -20513.072
$ws.Cells.Item(55, 8).Value = 991.05554  # LTW H55 was 1036.4706
$ws.Cells.Item(55, 10).Value = 862.8570999999999  # LTW J55 was 970.1667
$ws.Cells.Item(55, 12).Value = 862.8570999999999  # LTW L55 was 970.1667
$ws.Cells.Item(55, 14).Value = -1208.8571  # LTW N55 was -1316.1667

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 40388  # WVR H70 was 40414.285
$ws.Cells.Item(70, 10).Value = 40388  # WVR J70 was 40414.285
$ws.Cells.Item(70, 12).Value = 40388  # WVR L70 was 40414.285
$ws.Cells.Item(70, 14).Value = -41018  # WVR N70 was -41044.285
$ws.Cells.Item(73, 8).Value = 40388  # WVR H73 was 40414.285
$ws.Cells.Item(73, 10).Value = 40388  # WVR J73 was 40414.285
$ws.Cells.Item(73, 12).Value = 40388  # WVR L73 was 40414.285
$ws.Cells.Item(73, 14).Value = -42572  # WVR N73 was -42598.285
$ws.Cells.Item(122, 8).Value = 5115.1333  # WVR H122 was 5387.643
$ws.Cells.Item(122, 9).Value = 3658.6667  # WVR I122 was 3953.5
$ws.Cells.Item(122, 11).Value = 10976.0001  # WVR K122 was 11860.5
$ws.Cells.Item(122, 13).Value = -8526.000100000001  # WVR M122 was -9410.5
$ws.Cells.Item(136, 8).Value = 2931.3428  # WVR H136 was 3133.0312
$ws.Cells.Item(136, 9).Value = 2105.4138  # WVR I136 was 2152.0356
$ws.Cells.Item(136, 10).Value = 6923.3335  # WVR J136 was 10000
$ws.Cells.Item(136, 11).Value = 6316.241399999999  # WVR K136 was 6456.1068
$ws.Cells.Item(136, 12).Value = 20770.0005  # WVR L136 was 30000
$ws.Cells.Item(136, 13).Value = -3766.241399999999  # WVR M136 was -3906.1068
$ws.Cells.Item(136, 14).Value = -25870.0005  # WVR N136 was -35100
